$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.709.09'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").Value = '2.127.46'
$ws.Range("E3").Value = '  +10.86%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.75'
$ws.Range("E5").Value = '  +2.71%  '

$ws.Range("E6").Value = '  -4.18%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.15'
$ws.Range("E8").Value = '  +6.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.76'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("E10").Value = '  +2.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0744'
$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").Value = '2.439.85'
$ws.Range("E13").Value = '  +11.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.38'
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("E15").Value = '  +5.44%  '

$ws.Range("D16").Value = '2.129.80'
$ws.Range("E16").Value = '  +10.88%  '

$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").Value = '36.787.78'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.71'
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("D20").Value = '0.0₃0839'
$ws.Range("E20").Value = '  -2.71%  '

$ws.Range("E21").Value = '  +0.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '241.62'
$ws.Range("E22").Value = '  -4.05%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  +0.09%  '

$ws.Range("E25").Value = '  -7.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.84'
$ws.Range("E26").Value = '  +2.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.80'
$ws.Range("E27").Value = '  +16.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.26'
$ws.Range("E28").Value = '  +5.01%  '

$ws.Range("E29").Value = '  -7.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.82'
$ws.Range("E30").Value = '  +61.35%  '

$ws.Range("E31").Value = '  -4.28%  '

$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("E33").Value = '  +14.32%  '

$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("E35").Value = '  +17.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.954'
$ws.Range("E36").Value = '  +8.88%  '

$ws.Range("E37").Value = '  -4.65%  '

$ws.Range("E38").Value = '  -0.14%  '

$ws.Range("E39").Value = '  -4.29%  '

$ws.Range("E40").Value = '  -9.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.21'
$ws.Range("E41").Value = '  +8.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0226'
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.27'
$ws.Range("E43").Value = '  -7.01%  '

$ws.Range("E44").Value = '  +11.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.18'
$ws.Range("E45").Value = '  -5.85%  '

$ws.Range("D46").Value = '1.358.49'
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.26'
$ws.Range("E47").Value = '  +12.56%  '

$ws.Range("E48").Value = '  +3.51%  '

$ws.Range("D49").Value = '2.328.21'
$ws.Range("E49").Value = '  +11.07%  '

$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.30'
$ws.Range("E51").Value = '  -3.03%  '
